$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 8496
$ws1.Range("F6").Value = 310
$ws1.Range("F9").Value = 122
$ws1.Range("F12").Value = 895
$ws1.Range("F13").Value = 3673
$ws1.Range("F14").Value = 260
$ws1.Range("F15").Value = 145
$ws1.Range("F16").Value = 790
$ws1.Range("F17").Value = 773
$ws1.Range("F19").Value = 493
$ws1.Range("F22").Value = 1286
$ws1.Range("F24").Value = 430
$ws1.Range("F28").Value = 344
$ws1.Range("F33").Value = 653
$ws1.Range("F34").Value = 44
$ws1.Range("F36").Value = 76
$ws1.Range("F39").Value = 154
$ws1.Range("F40").Value = 5

$ws4.Range("F4").Value = 8496
$ws4.Range("F7").Value = 310
$ws4.Range("F10").Value = 122
$ws4.Range("F13").Value = 895
$ws4.Range("F15").Value = 3673
$ws4.Range("F16").Value = 260
$ws4.Range("F17").Value = 145
$ws4.Range("F19").Value = 790
$ws4.Range("F20").Value = 773
$ws4.Range("F23").Value = 493
$ws4.Range("F27").Value = 1286
$ws4.Range("F29").Value = 430
$ws4.Range("F34").Value = 344
$ws4.Range("F39").Value = 653
$ws4.Range("F40").Value = 44
$ws4.Range("F42").Value = 76
$ws4.Range("F45").Value = 154
$ws4.Range("F46").Value = 5
